$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

$ws.Cells.Item(2, 6).Value = "KN99alpha"
$ws.Cells.Item(5, 6).Value = "TDY2258"
$ws.Cells.Item(8, 6).Value = "TYS2271"

$ws.Range("F11").Select()
